$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A26").Value = 111439199
$ws.Range("B26").Value = 77515
$ws.Range("E26").Value = 6425
$ws.Range("F26").Value = "Garnlav"
$ws.Range("G26").Value = "Alectoria sarmentosa"
$ws.Range("H26").Value = "(Ach.) Ach."
$ws.Range("Q26").Value = 571460.5263831713
$ws.Range("R26").Value = 6988072.746377321

$ws.Range("A27").Value = 111440132
$ws.Range("B27").Value = 77515
$ws.Range("E27").Value = 6425
$ws.Range("F27").Value = "Garnlav"
$ws.Range("G27").Value = "Alectoria sarmentosa"
$ws.Range("H27").Value = "(Ach.) Ach."
$ws.Range("Q27").Value = 571561.4085783331
$ws.Range("R27").Value = 6987995.007176694

$ws.Range("A28").Value = 111438212
$ws.Range("B28").Value = 78578
$ws.Range("E28").Value = 6458
$ws.Range("F28").Value = "Lunglav"
$ws.Range("G28").Value = "Lobaria pulmonaria"
$ws.Range("H28").Value = "(L.) Hoffm."
$ws.Range("Q28").Value = 571110.4113759904
$ws.Range("R28").Value = 6988002.384436167

$ws.Range("A29").Value = 111439399
$ws.Range("B29").Value = 78578
$ws.Range("E29").Value = 6458
$ws.Range("F29").Value = "Lunglav"
$ws.Range("G29").Value = "Lobaria pulmonaria"
$ws.Range("H29").Value = "(L.) Hoffm."
$ws.Range("Q29").Value = 571513.3871347769
$ws.Range("R29").Value = 6988089.353432911

$ws.Range("A30").Value = 111440415
$ws.Range("B30").Value = 78578
$ws.Range("E30").Value = 6458
$ws.Range("F30").Value = "Lunglav"
$ws.Range("G30").Value = "Lobaria pulmonaria"
$ws.Range("H30").Value = "(L.) Hoffm."
$ws.Range("Q30").Value = 571485.2187846474
$ws.Range("R30").Value = 6987922.006598755

$ws.Range("A31").Value = 111439000
$ws.Range("B31").Value = 77515
$ws.Range("E31").Value = 6425
$ws.Range("F31").Value = "Garnlav"
$ws.Range("G31").Value = "Alectoria sarmentosa"
$ws.Range("H31").Value = "(Ach.) Ach."
$ws.Range("Q31").Value = 571401.1370697253
$ws.Range("R31").Value = 6988021.924799537

$ws.Range("A32").Value = 111439556
$ws.Range("B32").Value = 78578
$ws.Range("E32").Value = 6458
$ws.Range("F32").Value = "Lunglav"
$ws.Range("G32").Value = "Lobaria pulmonaria"
$ws.Range("H32").Value = "(L.) Hoffm."
$ws.Range("Q32").Value = 571541.5030949884
$ws.Range("R32").Value = 6988072.708431445

$ws.Range("A33").Value = 111440237
$ws.Range("B33").Value = 77515
$ws.Range("E33").Value = 6425
$ws.Range("F33").Value = "Garnlav"
$ws.Range("G33").Value = "Alectoria sarmentosa"
$ws.Range("H33").Value = "(Ach.) Ach."
$ws.Range("Q33").Value = 571532.8559679713
$ws.Range("R33").Value = 6987969.392762338

$ws.Range("A34").Value = 111438249
$ws.Range("B34").Value = 77267
$ws.Range("E34").Value = 6446
$ws.Range("F34").Value = "Kolflarnlav"
$ws.Range("G34").Value = "Carbonicola anthracophila"
$ws.Range("H34").Value = "(Nyl.) Bendiksby & Timdal"
$ws.Range("Q34").Value = 571110.4113759904
$ws.Range("R34").Value = 6988002.384436167

$ws.Range("A35").Value = 111438897
$ws.Range("B35").Value = 78578
$ws.Range("E35").Value = 6458
$ws.Range("F35").Value = "Lunglav"
$ws.Range("G35").Value = "Lobaria pulmonaria"
$ws.Range("H35").Value = "(L.) Hoffm."
$ws.Range("Q35").Value = 571342.05486291
$ws.Range("R35").Value = 6988019.26613469

$ws.Range("A36").Value = 111438342
$ws.Range("B36").Value = 77268
$ws.Range("E36").Value = 228912
$ws.Range("F36").Value = "Mörk kolflarnlav"
$ws.Range("G36").Value = "Carbonicola myrmecina"
$ws.Range("H36").Value = "(Ach.) Bendiksby & Timdal"
$ws.Range("Q36").Value = 571171.3525270248
$ws.Range("R36").Value = 6988003.262354148

$ws.Range("A37").Value = 111438823
$ws.Range("B37").Value = 78578
$ws.Range("E37").Value = 6458
$ws.Range("F37").Value = "Lunglav"
$ws.Range("G37").Value = "Lobaria pulmonaria"
$ws.Range("H37").Value = "(L.) Hoffm."
$ws.Range("Q37").Value = 571355.5533441483
$ws.Range("R37").Value = 6987984.581054003

$ws.Range("A38").Value = 111438961
$ws.Range("B38").Value = 89686
$ws.Range("E38").Value = 658
$ws.Range("F38").Value = "Rosenticka"
$ws.Range("G38").Value = "Rhodofomes roseus"
$ws.Range("H38").Value = "(Alb. & Schwein.) Kotl. & Pouzar"
$ws.Range("Q38").Value = 571355.5533441483
$ws.Range("R38").Value = 6987984.581054003

$ws.Range("A39").Value = 111440195
$ws.Range("B39").Value = 78578
$ws.Range("E39").Value = 6458
$ws.Range("F39").Value = "Lunglav"
$ws.Range("G39").Value = "Lobaria pulmonaria"
$ws.Range("H39").Value = "(L.) Hoffm."
$ws.Range("Q39").Value = 571540.7114327301
$ws.Range("R39").Value = 6987984.557469011

$ws.Range("A40").Value = 111440357
$ws.Range("B40").Value = 78578
$ws.Range("E40").Value = 6458
$ws.Range("F40").Value = "Lunglav"
$ws.Range("G40").Value = "Lobaria pulmonaria"
$ws.Range("H40").Value = "(L.) Hoffm."
$ws.Range("Q40").Value = 571472.8873933584
$ws.Range("R40").Value = 6987944.905536451
